$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.676.56'
$ws.Range('E2').Value = '  -1.82%  '
$ws.Range('D3').Value = '2.340.04'
$ws.Range('E3').Value = '  -2.00%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'502.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').Value = "'128.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -2.63%  '
$ws.Range('D9').Value = '2.346.13'
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('D10').Value = "'0.0970"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.13%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').Value = "'4.76"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.85%  '
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').Value = '2.756.74'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').Value = '55.627.26'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').Value = "'21.58"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '2.355.63'
$ws.Range('E18').Value = '  -2.86%  '
$ws.Range('D19').Value = "'9.90"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.78%  '
$ws.Range('D20').Value = "'309.92"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.46%  '
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').Value = "'6.22"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.61%  '
$ws.Range('D23').Value = "'0.997"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').Value = "'65.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('D25').Value = "'0.998"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = "'0.369"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').Value = "'0.145"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('E28').Value = '  -4.55%  '
$ws.Range('D29').Value = "'171.17"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').Value = '0.0₃0700'
$ws.Range('E31').Value = '  -3.03%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('D33').Value = "'5.74"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.01%  '
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -5.55%  '
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').Value = "'1.16"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('D38').Value = "'3.62"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.52%  '
$ws.Range('D39').Value = "'0.819"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').Value = "'36.02"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.11%  '
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('E42').Value = '  -1.00%  '
$ws.Range('D43').Value = "'125.91"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.09%  '
$ws.Range('E44').Value = '  -3.49%  '
$ws.Range('D45').Value = "'0.553"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.14%  '
$ws.Range('D46').Value = "'0.0889"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').Value = "'236.39"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.70%  '
$ws.Range('D48').Value = "'0.0473"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.23%  '
$ws.Range('E49').Value = '  -2.13%  '
$ws.Range('D50').Value = "'16.68"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = "'0.952"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '
